$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SVM)
$ws.Range("B2").Value = 0.667
$ws.Range("D2").Value = 0.741

# Row 3 (LR)
$ws.Range("B3").Value = 0.6820000000000001
$ws.Range("C3").Value = 0.002
$ws.Range("D3").Value = 0.744
$ws.Range("F3").Value = 0.002020220831301504

# Row 4 (LDA)
$ws.Range("B4").Value = 0.643
$ws.Range("D4").Value = 0.71

# Row 5 (RF)
$ws.Range("B5").Value = 0.5580000000000001
$ws.Range("C5").Value = 0.026
$ws.Range("D5").Value = 0.573
$ws.Range("E5").Value = 0.018
$ws.Range("F5").Value = 0.02630944998180856

# Row 6 (AB)
$ws.Range("B6").Value = 0.476
$ws.Range("D6").Value = 0.466

# Row 7 (KNN)
$ws.Range("B7").Value = 0.635
$ws.Range("D7").Value = 0.681

# Row 8 (GNB)
$ws.Range("B8").Value = 0.643
$ws.Range("D8").Value = 0.722
$ws.Range("F8").Value = 0
